# Actualización automática 2025-06-18 11:45:09
#
# Updates the sales figures for client "PADILLA MIER BERTHA MARIETA"
# (advisor HIDALGO HIDALGO PEDRO GUSTAVO) to reflect a new PORCELANATO
# sale of 45.91 registered in June, and propagates the corresponding
# totals / ratios across the three report sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet "VENTAS POR GRUPO" ---------------------------------------------
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsGrupo.Range("M16").Value = 45.91
$wsGrupo.Range("M22").Value = "1 de 20"

# --- Sheet "VENTA MENSUAL" -------------------------------------------------
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsMensual.Range("F16").Value = 45.91
$wsMensual.Range("F22").Value = -8.849999999999994

# --- Sheet "CUMPLIMIENTO MENSUAL" -----------------------------------------
$wsCumplimiento = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$wsCumplimiento.Range("D16").Value = 45.91
$wsCumplimiento.Range("E16").Value = 29486.53
$wsCumplimiento.Range("F16").Value = 0.001554561695545644

$wsCumplimiento.Range("D19").Value = -8.849999999999966
$wsCumplimiento.Range("E19").Value = 50396.04762291769
$wsCumplimiento.Range("F19").Value = -0.0001756398533260502
